{"js": "// Add a new bullet (\"technologies used: ...\") right after the paragraph\n// ending in \"... science students about software development in the real\n// world\", as a sibling bullet (same Compact style, same numbering list\n// numId=1005, ilvl=1) in the \"software architect\" job entry.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its (unique) text content.\nconst anchorText =\n  \"science students about software development in the real world\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph for insertion\");\n}\n\n// Read the anchor's list membership so the new bullet matches it exactly.\nanchor.load(\"style\");\nanchor.load(\"listItemOrNullObject/level\");\nawait context.sync();\n\nconst anchorStyle = anchor.style;\nconst listItem = anchor.listItemOrNullObject;\nconst level = listItem.isNullObject ? 1 : listItem.level;\n\nconst list = anchor.listOrNullObject;\nlist.load(\"id\");\nawait context.sync();\nconst listId = list.id;\n\n// Insert the new paragraph right after the anchor, copy its style, and\n// attach it to the same numbered list at the same indentation level.\nconst newParagraph = anchor.insertParagraph(\n  \"technologies used: C#, WPF, ASP.NET MVC, SQL Server, NHibernate, SignalR etc.\",\n  \"After\"\n);\nnewParagraph.style = anchorStyle;\nnewParagraph.attachToList(listId, level);\n\nawait context.sync();\n", "ps1": "# Add a new bullet (\"technologies used: ...\") right after the paragraph\n# ending in \"... science students about software development in the real\n# world\", as a sibling bullet in the same numbered list (same style and\n# same numbering level) within the \"software architect\" job entry.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"science students about software development in the real world\"\n\n# Use Find on a Range to locate the anchor paragraph without relying on a\n# hard-coded paragraph index.\n$findRange = $d.Content\n$found = $findRange.Find.Execute($anchorText)\nif (-not $found) {\n    throw \"Could not find anchor paragraph for insertion\"\n}\n\n# Resolve the matching Paragraph object by locating the paragraph whose\n# range contains the start of the found text.\n$anchor = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($findRange.Start -ge $candidate.Range.Start -and $findRange.Start -lt $candidate.Range.End) {\n        $anchor = $candidate\n        break\n    }\n}\nif (-not $anchor) {\n    throw \"Could not resolve anchor paragraph object\"\n}\n\n# Insert a new paragraph right after the anchor. InsertParagraphAfter\n# copies the anchor's paragraph formatting (style + numbering), so the\n# new bullet automatically lands in the same list/level.\n$anchor.Range.InsertParagraphAfter()\n\n$newPara = $anchor.Next()\n$newPara.Range.Text = \"technologies used: C#, WPF, ASP.NET MVC, SQL Server, NHibernate, SignalR etc.\"\n"}
